# Saldo.xlsx update
# - Merge the "Patricia" (004421636) and "Cairo" (005905737) rows into a
#   single updated "Andre" (004514241) row with balance 22246.37
# - Update "Pedro" (004460487) balance from 6970.94 to 15000
# - Remove the old duplicate "Andre" (004514241) row that had balance 1062.55

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update Pedro's balance (row 6) while row indices are still in their
#    original positions.
$ws.Cells.Item(6, 3).Value = 15000

# 2) Turn the Patricia row (row 4) into the consolidated Andre row.
#    The leading apostrophe forces Excel to keep the account number as
#    text (preserving the leading zeros) instead of coercing it to a
#    number.
$ws.Cells.Item(4, 1).Value = "'004514241"
$ws.Cells.Item(4, 2).Value = "Andre"
$ws.Cells.Item(4, 3).Value = 22246.37

# 3) Delete the old duplicate Andre row (row 10) first since it is below
#    every row touched above, so those rows keep their indices valid.
$ws.Rows(10).Delete()

# 4) Delete the now-redundant Cairo row (row 5); Patricia's row already
#    carries the merged Andre data, and everything below shifts up.
$ws.Rows(5).Delete()

Write-Host "Saldo.xlsx update applied"
